$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table runs from row 2 through the last used row; find it dynamically
# (falls back to the known last row of data, 211, if something odd happens)
# and append one new daily record right below it, matching the existing
# column layout (A:DX).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 211 }
$newRow = $lastRow + 1

# New record for 2020-10-02 (Excel serial date 44106), one value per column
# A:DX, in the same order as the existing header row.
$rowValues = @(44106,841531,2744,116663,67748,271929,29305,7472,6144,8598,9724,20511,4023,24187,33840,8230,11966,15470,15496,18602,16172,3822,3974,11294,31727,14127,12592,63784,2512,1497,820,479,933,533,843,2078,5999,38287,10351,2579,49096,1173,23399,1537,10674,1686,1617,8786,2069,966,2503,2701,67508,14318,7197,10073,7574,257,1468,2757,745,2180,10069,9607,10836,14404,1976,912,14441,11807,13820,3391,2367,6360,5121,2795,6301,3963,2364,1310,3099,2271,2163,2071,6688,2285,1544,1899,2187,2345,2766,1885,1231,1237,1156,3494,1542,999,1192,1775,1652,823,927,1410,1857,1742,1713,1314,335,373,841,792,511,544,392,699,765,534,514,374,527,141264,354929,21685,154967,96133,47948,13353)

# Marshal the PowerShell array into a (1 x N) COM-friendly 2-D array and
# write it to the new row in a single call so the sheet's per-column default
# styles (date format in column A, the shared style on BS:DQ, etc.) are
# picked up automatically from the column definitions, just like typing the
# values in by hand would do.
$colCount = $rowValues.Length
$arr = New-Object 'object[,]' 1,$colCount
for ($i = 0; $i -lt $colCount; $i++) {
    $arr[0, $i] = $rowValues[$i]
}

$startCell = $ws.Cells.Item($newRow, 1)
$endCell = $ws.Cells.Item($newRow, $colCount)
$dstRange = $ws.Range($startCell, $endCell)
$dstRange.Value = $arr

# Move the active selection to where the author left off after entering the
# new row.
$ws.Range("CI202").Select() | Out-Null
